# Fixed bugs in producing trees and mapping process.
# The updated data adds one more row (row 6) of results to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row
$ws.Range("A6").Value = 0.64453125
$ws.Range("B6").Value = 0.68619793653488104

# Scroll the window so row 4 is the first visible row, matching the
# author's on-screen state when the file was saved.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# Move the active selection to the newly entered cell
$ws.Range("B6").Select()
